$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in A2 and A3 (city names)
$ws.Range("A2").Value = "Delhi"
$ws.Range("A3").Value = "Dubai"

# Update the selection to A2
$ws.Range("A2").Select()
